# Weekly CompStat data refresh: bump volume/week numbers and update the
# crime statistics table (rows 15-30) with newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text updates (rich text cells) - edit just the changed runs
# ---------------------------------------------------------------------

# "Volume 32   Number  34" -> "...Number  35"
$ws.Range("A8").Characters(21, 2).Text = "35"

# "Report Covering the Week  8/18/2025  Through  8/24/2025"
#   -> "...8/25/2025  Through  8/31/2025"
$ws.Range("C9").Characters(27, 9).Text = "8/25/2025"
$ws.Range("C9").Characters(47, 9).Text = "8/31/2025"

# ---------------------------------------------------------------------
# Crime statistics table updates (rows 15-30)
# ---------------------------------------------------------------------

# Row 15
$ws.Range("N15").Value = -45.454545454545

# Row 16
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 33.333333333333
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -18.181818181818
$ws.Range("I16").Value = 91
$ws.Range("J16").Value = 115
$ws.Range("K16").Value = -20.869565217391
$ws.Range("L16").Value = -21.551724137931
$ws.Range("M16").Value = -51.336898395721
$ws.Range("N16").Value = -88.135593220339

# Row 17
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -9.090909090909
$ws.Range("I17").Value = 206
$ws.Range("J17").Value = 230
$ws.Range("K17").Value = -10.434782608695
$ws.Range("L17").Value = -22.556390977443
$ws.Range("M17").Value = 28.75
$ws.Range("N17").Value = -8.849557522123

# Row 18
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 8
$ws.Range("H18").Value = 33.333333333333
$ws.Range("I18").Value = 68
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = 13.333333333333
$ws.Range("L18").Value = -31.313131313131
$ws.Range("M18").Value = -67.924528301886
$ws.Range("N18").Value = -92.804232804232

# Row 19
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = -14.285714285714
$ws.Range("I19").Value = 209
$ws.Range("J19").Value = 230
$ws.Range("K19").Value = -9.130434782608
$ws.Range("L19").Value = -4.128440366972
$ws.Range("M19").Value = -4.128440366972
$ws.Range("N19").Value = -43.817204301075

# Row 20
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = 42.857142857142
$ws.Range("F20").Value = 33
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = 65
$ws.Range("I20").Value = 221
$ws.Range("J20").Value = 192
$ws.Range("K20").Value = 15.104166666666
$ws.Range("L20").Value = 15.706806282722
$ws.Range("M20").Value = 25.568181818181
$ws.Range("N20").Value = -90.340909090909

# Row 21 (Total Seven Major Felony row)
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -20
$ws.Range("F21").Value = 95
$ws.Range("G21").Value = 88
$ws.Range("H21").Value = 7.954545454545
$ws.Range("I21").Value = 810
$ws.Range("J21").Value = 847
$ws.Range("K21").Value = -4.368358913813
$ws.Range("L21").Value = -10.596026490066
$ws.Range("M21").Value = -16.149068322981
$ws.Range("N21").Value = -82.482698961937

# Row 22 (Transit) - several cells flip between numeric and the "0"/"***.*"
# placeholder text used elsewhere on this sheet.
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("C22").Value = 1

$ws.Range("G22").Value = "0"
$ws.Range("G22").NumberFormat = "General"

$ws.Range("H22").Value = "***.*"
$ws.Range("H22").NumberFormat = "General"

$ws.Range("I22").Value = 10
$ws.Range("K22").Value = 66.666666666666
$ws.Range("L22").Value = 25
$ws.Range("M22").Value = -9.090909090909

# Row 24
$ws.Range("C24").Value = 40
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 33.333333333333
$ws.Range("F24").Value = 108
$ws.Range("G24").Value = 97
$ws.Range("H24").Value = 11.340206185567
$ws.Range("I24").Value = 840
$ws.Range("J24").Value = 920
$ws.Range("K24").Value = -8.695652173913
$ws.Range("L24").Value = -15.237134207870
$ws.Range("M24").Value = 63.742690058479

# Row 25
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = 37.5
$ws.Range("F25").Value = 47
$ws.Range("G25").Value = 46
$ws.Range("H25").Value = 2.173913043478
$ws.Range("I25").Value = 402
$ws.Range("J25").Value = 431
$ws.Range("K25").Value = -6.728538283062
$ws.Range("L25").Value = 19.642857142857

# Row 26
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = 55.555555555555
$ws.Range("F26").Value = 33
$ws.Range("G26").Value = 43
$ws.Range("H26").Value = -23.255813953488
$ws.Range("I26").Value = 367
$ws.Range("J26").Value = 424
$ws.Range("K26").Value = -13.443396226415
$ws.Range("L26").Value = 0.273224043715
$ws.Range("M26").Value = -13.033175355450

# Row 27 (Misd. Assault) - D27/E27 flip from numeric to placeholder text
$ws.Range("D27").Value = "0"
$ws.Range("D27").NumberFormat = "General"

$ws.Range("E27").Value = "***.*"
$ws.Range("E27").NumberFormat = "General"

$ws.Range("L27").Value = -16

# Row 28
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = -71.428571428571
$ws.Range("L28").Value = 13.953488372093

# Row 29
$ws.Range("N29").Value = -86.363636363636

# Row 30
$ws.Range("N30").Value = -83.333333333333
